$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.494470477104187
$ws.Range("B1").Value = 1.727225422859192
$ws.Range("C1").Value = 2.186141967773438
$ws.Range("D1").Value = 3.510639905929565
$ws.Range("E1").Value = 3.901979207992554
